$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- CAD BOM sheet updates ----
# Record the purchase quantity ("1 (pack of 100)") for the screw and nut rows.
$ws.Range("D3").Value = "1 (pack of 100)"
$ws.Range("D4").Value = "1 (pack of 100)"

# Turn the plain-text McMaster-Carr URLs in the Link column into real hyperlinks.
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.mcmaster.com/90128A179/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.mcmaster.com/90591A111/") | Out-Null

# Resize column D closer to column E's width now that it holds real content.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(5).ColumnWidth()

# Restore the cursor/selection position recorded at save time.
$ws.Range("N13").Select() | Out-Null
